# Crisp DM lecture notes added, slide 1 finished
# Fill in the answers on the "Table -1.1" sheet (first worksheet) for the
# "Understand the Data Set" question block (rows 5-9, column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table -1.1")

$ws.Range("C5").Value = 66370
$ws.Range("C6").Value = 66368
$ws.Range("C7").Value = "permalink"
$ws.Range("C8").Value = "Y"
$ws.Range("C9").Value = 114954

# Make this sheet the active / selected tab, matching the author's final
# view state (was on "Table-5.1" before, now back on "Table -1.1").
$ws.Select()
$ws.Range("B21").Select()
